# Commit: "added examples for all profiles, edited medication guidance"
# 1. Update the Date property on the Metadata sheet.
# 2. On the Concepts sheet, strip the "$CADSR:" prefix from each Code cell
#    (keeping them as text, not numbers) and fill in the previously-empty
#    Definition column with example/definition text.

$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update Date value (B8) ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2022-09-02T15:43:08-05:00"

# --- Concepts sheet: strip "$CADSR:" prefix from codes, add definitions ---
$concepts = $wb.Worksheets.Item("Concepts")

# Force the Code column to stay text (these codes are numeric-looking, but
# are identifiers, not numbers) before writing the stripped-down values.
$codes = $concepts.Range("B2:B6")
$codes.NumberFormat = "@"

$concepts.Range("B2").Value = "4722619"
$concepts.Range("D2").Value = "Acute Graft Vs Host Disease"

$concepts.Range("B3").Value = "2991794"
$concepts.Range("D3").Value = "Classic chronic Graft Vs Host Disease"

$concepts.Range("B4").Value = "2991795"
$concepts.Range("D4").Value = "Overlap chronic Graft Vs Host Disease"

$concepts.Range("B5").Value = "2991796"
$concepts.Range("D5").Value = "Late Acute Graft Vs Host Disease"

$concepts.Range("B6").Value = "3261243"
$concepts.Range("D6").Value = "Not Applicable"
